$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "8Winter2023"
$ws.Range("C5").Value = "9Winter2023"
$ws.Range("C6").Value = "10Winter2023"
$ws.Range("C7").Value = "11WInter2023"
